$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws "D2" "43.151.02"
$ws.Range("E2").Value = "  +0.70%  "

Set-TextCell $ws "D3" "2.387.41"
$ws.Range("E3").Value = "  +6.05%  "

$ws.Range("E4").Value = "  -0.54%  "

Set-TextCell $ws "D5" "326.61"
$ws.Range("E5").Value = "  +10.64%  "

Set-TextCell $ws "D6" "105.88"
$ws.Range("E6").Value = "  -6.96%  "

Set-TextCell $ws "D7" "0.645"
$ws.Range("E7").Value = "  +2.62%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("E9").Value = "  +8.56%  "

Set-TextCell $ws "D10" "41.88"
$ws.Range("E10").Value = "  -5.47%  "

$ws.Range("E11").Value = "  +1.61%  "

Set-TextCell $ws "D12" "8.67"
$ws.Range("E12").Value = "  -3.17%  "

Set-TextCell $ws "D13" "1.05"
$ws.Range("E13").Value = "  -1.74%  "

Set-TextCell $ws "D14" "17.22"
$ws.Range("E14").Value = "  +13.89%  "

$ws.Range("E15").Value = "  +2.13%  "

Set-TextCell $ws "D16" "2.748.23"
$ws.Range("E16").Value = "  +6.20%  "

Set-TextCell $ws "D17" "2.389.31"
$ws.Range("E17").Value = "  +6.38%  "

Set-TextCell $ws "D18" "43.056.38"
$ws.Range("E18").Value = "  +0.56%  "

Set-TextCell $ws "D19" "7.92"
$ws.Range("E19").Value = "  +9.35%  "

$ws.Range("E20").Value = "  +2.36%  "

Set-TextCell $ws "D21" "76.88"
$ws.Range("E21").Value = "  +3.14%  "

Set-TextCell $ws "D22" "3.62"
$ws.Range("E22").Value = "  +7.07%  "

Set-TextCell $ws "D23" "278.70"
$ws.Range("E23").Value = "  +11.33%  "

$ws.Range("E24").Value = "  +0.18%  "

Set-TextCell $ws "D25" "9.61"
$ws.Range("E25").Value = "  +6.86%  "

$ws.Range("E26").Value = "  +2.04%  "

$ws.Range("E27").Value = "  -0.02%  "

Set-TextCell $ws "D28" "23.23"
$ws.Range("E28").Value = "  +5.65%  "

Set-TextCell $ws "D29" "37.92"
$ws.Range("E29").Value = "  +0.41%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws "D30" "175.18"
$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws "D31" "2.18"
$ws.Range("E31").Value = "  -1.58%  "

Set-TextCell $ws "D32" "3.18"
$ws.Range("E32").Value = "  +1.67%  "

$ws.Range("E33").Value = "  +5.21%  "

Set-TextCell $ws "D34" "5.89"
$ws.Range("E34").Value = "  +3.30%  "

$ws.Range("E35").Value = "  +5.35%  "

$ws.Range("E36").Value = "  -3.20%  "

Set-TextCell $ws "D37" "4.23"
$ws.Range("E37").Value = "  -1.05%  "

$ws.Range("E38").Value = "  -1.93%  "

Set-TextCell $ws "D39" "0.108"
$ws.Range("E39").Value = "  +2.89%  "

$ws.Range("E40").Value = "  +15.98%  "

$ws.Range("E41").Value = "  +19.47%  "

$ws.Range("E42").Value = "  +1.55%  "

Set-TextCell $ws "D43" "69.90"
$ws.Range("E43").Value = "  -3.43%  "

Set-TextCell $ws "D44" "122.51"
$ws.Range("E44").Value = "  +15.95%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell $ws "D45" "1.00"
$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextCell $ws "D46" "92.94"
$ws.Range("E46").Value = "  +58.52%  "

Set-TextCell $ws "D47" "12.50"
$ws.Range("E47").Value = "  -0.28%  "

Set-TextCell $ws "D48" "9.40"
$ws.Range("E48").Value = "  +9.02%  "

Set-TextCell $ws "D49" "5.53"
$ws.Range("E49").Value = "  +0.20%  "

Set-TextCell $ws "D50" "1.32"
$ws.Range("E50").Value = "  +0.96%  "

Set-TextCell $ws "D51" "1.596.35"
